$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shift readings logged after the previous batch (rows 5-8).
# Columns: A = shift timestamp, B = whole tomatoes, C = total tomatoes.
$newRows = @(
    @{ Row = 5; Time = "11/18/2022 20:50:36"; Whole = 107; Total = 107 },
    @{ Row = 6; Time = "11/18/2022 21:20:03"; Whole = 150; Total = 150 },
    @{ Row = 7; Time = "11/18/2022 21:24:51"; Whole = 110; Total = 110 },
    @{ Row = 8; Time = "11/18/2022 21:27:11"; Whole = 115; Total = 115 }
)

# Tomato amounts at/above this threshold are flagged as too large.
$tooLargeThreshold = 110

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Range("A$r").Value = $entry.Time
    $ws.Range("B$r").Value = $entry.Whole
    $ws.Range("C$r").Value = $entry.Total

    if ($entry.Whole -ge $tooLargeThreshold) {
        # Color the whole-tomatoes cell red when the amount is too large.
        $ws.Range("B$r").Interior.Color = 255
    }
}
